$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @(6160,6080,6010,5940,5900,5870,5840,5830,5810,5810,5810,5810,5800,5790,5790,5780,5790,5800,5830,5880,5950,6030,6120,6220,6320,6410,6500,6570,6620,6660,6680,6680,6660,6640,6610,6570,6530,6500,6470,6450,6440,6430,6430,6430,6440,6460,6470,6500,6520,6550,6580,6610,6640,6670,6700,6730,6750,6780,6810,6850,6900,6970,7040,7120,7200,7280,7360,7440,7510,7600,7690,7780,7860,7920,7980,8010,8040,8080,8100,8120,8100,8070,8030,7970,7900,7800,7650,7480,7300,7110,6980,6860,6490,6450,6390,6290)
$bValues = @(45863,45863.01041666666,45863.02083333334,45863.03125,45863.04166666666,45863.05208333334,45863.0625,45863.07291666666,45863.08333333334,45863.09375,45863.10416666666,45863.11458333334,45863.125,45863.13541666666,45863.14583333334,45863.15625,45863.16666666666,45863.17708333334,45863.1875,45863.19791666666,45863.20833333334,45863.21875,45863.22916666666,45863.23958333334,45863.25,45863.26041666666,45863.27083333334,45863.28125,45863.29166666666,45863.30208333334,45863.3125,45863.32291666666,45863.33333333334,45863.34375,45863.35416666666,45863.36458333334,45863.375,45863.38541666666,45863.39583333334,45863.40625,45863.41666666666,45863.42708333334,45863.4375,45863.44791666666,45863.45833333334,45863.46875,45863.47916666666,45863.48958333334,45863.5,45863.51041666666,45863.52083333334,45863.53125,45863.54166666666,45863.55208333334,45863.5625,45863.57291666666,45863.58333333334,45863.59375,45863.60416666666,45863.61458333334,45863.625,45863.63541666666,45863.64583333334,45863.65625,45863.66666666666,45863.67708333334,45863.6875,45863.69791666666,45863.70833333334,45863.71875,45863.72916666666,45863.73958333334,45863.75,45863.76041666666,45863.77083333334,45863.78125,45863.79166666666,45863.80208333334,45863.8125,45863.82291666666,45863.83333333334,45863.84375,45863.85416666666,45863.86458333334,45863.875,45863.88541666666,45863.89583333334,45863.90625,45863.91666666666,45863.92708333334,45863.9375,45863.94791666666,45863.95833333334,45863.96875,45863.97916666666,45863.98958333334)

for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}
